$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A; this shifts all existing columns (A:W) to (B:X)
# and shifts the merged cell ranges / dimension accordingly.
$ws.Columns("A").Insert()

# Rows 1, 3 and 20 are hidden; temporarily unhide them so writing a value
# doesn't trigger an auto row-height adjustment, then re-hide them.
$ws.Rows(1).Hidden = $false
$ws.Rows(3).Hidden = $false
$ws.Rows(20).Hidden = $false

# New "Match ID" column: bold header/data font, no border (matches the
# existing bold-no-border style used elsewhere), value 25 for every
# data row (4-19), and an unstyled 25 on the hidden totals row (20).
$ws.Range("A2:A19").Font.Bold = $true

$ws.Range("A2").Value = "Match ID"
$ws.Range("A4:A19").Value = 25
$ws.Range("A20").Value = 25

# Re-hide the rows that were hidden originally.
$ws.Rows(1).Hidden = $true
$ws.Rows(3).Hidden = $true
$ws.Rows(20).Hidden = $true

# Restore the original selection (insert + edits above move it around).
[void]$ws.Range("A2:A19").Select()
